# Add "Wins", "Losses", "Ties" columns (AD, AE, AF) to the season-record sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — match the styling of the other header cells (bold,
# centered, top-aligned, thin box border).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows 2-57 — every team row shares the same season record.
$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
